$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.932.28"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.892.08"
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.8333"
$ws.Range("E5").Value = "  +8.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.66"
$ws.Range("E6").Value = "  +0.55%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3255"
$ws.Range("E8").Value = "  +6.77%  "
$ws.Range("E9").Value = "  +5.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07028"
$ws.Range("E10").Value = "  +2.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08022"
$ws.Range("E11").Value = "  +0.54%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7482"
$ws.Range("E12").Value = "  +1.44%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.898.13"
$ws.Range("E13").Value = "  +0.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.196"
$ws.Range("E14").Value = "  +0.52%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.27"
$ws.Range("E15").Value = "  +1.35%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.935.73"
$ws.Range("E16").Value = "  +0.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.05"
$ws.Range("E17").Value = "  +1.97%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.916"
$ws.Range("E18").Value = "  +0.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.50"
$ws.Range("E19").Value = "  -0.83%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007764"
$ws.Range("E20").Value = "  +0.87%  "
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.149.36"
$ws.Range("E22").Value = "  +0.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("E24").Value = "  +0.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1612"
$ws.Range("E25").Value = "  +25.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.87"
$ws.Range("E26").Value = "  +0.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.194"
$ws.Range("E27").Value = "  -0.65%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.85"
$ws.Range("E28").Value = "  +1.08%  "
$ws.Range("E29").Value = "  +2.26%  "
$ws.Range("E30").Value = "  -1.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.519"
$ws.Range("E31").Value = "  +0.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.261"
$ws.Range("E32").Value = "  -0.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05622"
$ws.Range("E33").Value = "  +6.75%  "
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.278"
$ws.Range("E35").Value = "  +2.46%  "
$ws.Range("E36").Value = "  +1.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.717"
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("E38").Value = "  +0.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.776"
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("E40").Value = "  +0.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "71.94"
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.935"
$ws.Range("E42").Value = "  -4.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8430"
$ws.Range("E43").Value = "  +1.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9997"
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.890"
$ws.Range("E45").Value = "  +0.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.29"
$ws.Range("E46").Value = "  +1.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.592"
$ws.Range("E47").Value = "  -0.15%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "993.71"
$ws.Range("E48").Value = "  +9.44%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.678"
$ws.Range("E49").Value = "  -0.98%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.047.41"
$ws.Range("E50").Value = "  +0.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "35.98"
$ws.Range("E51").Value = "  -0.35%  "
